$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "Time to Discuss"

# Fill column E (rows 2-16) with "Breaktime"
for ($i = 2; $i -le 16; $i++) {
    $ws.Range("E" + $i).Value = "Breaktime"
}

# Give column E its own (narrower) width, splitting it off from column D
$ws.Columns.Item(5).ColumnWidth = 16

# Move the active selection, matching the saved view state
$ws.Range("I8").Select()
